{"js": "// Replace each three-digit-divided-by-one-digit expression in the\n// document with its updated value, per the commit's regenerated\n// worksheet numbers. Every occurrence is a unique, literal string\n// (e.g. \"717\u00f79=\"), so a direct search + replace-in-place keeps the\n// existing run formatting (font, size, etc.) untouched.\nconst replacements = [\n  [\"717\u00f79=\", \"147\u00f76=\"],\n  [\"632\u00f77=\", \"914\u00f75=\"],\n  [\"169\u00f74=\", \"873\u00f77=\"],\n  [\"852\u00f79=\", \"153\u00f72=\"],\n  [\"473\u00f78=\", \"562\u00f72=\"],\n  [\"531\u00f74=\", \"235\u00f74=\"],\n  [\"284\u00f72=\", \"795\u00f73=\"],\n  [\"431\u00f79=\", \"503\u00f72=\"],\n  [\"946\u00f78=\", \"891\u00f77=\"],\n  [\"804\u00f79=\", \"267\u00f76=\"],\n  [\"976\u00f76=\", \"568\u00f77=\"],\n  [\"951\u00f78=\", \"676\u00f73=\"],\n  [\"169\u00f72=\", \"405\u00f78=\"],\n  [\"904\u00f78=\", \"583\u00f75=\"],\n  [\"694\u00f72=\", \"869\u00f73=\"],\n  [\"863\u00f75=\", \"934\u00f72=\"],\n  [\"724\u00f72=\", \"542\u00f75=\"],\n  [\"116\u00f78=\", \"814\u00f74=\"],\n  [\"861\u00f73=\", \"763\u00f77=\"],\n  [\"217\u00f73=\", \"840\u00f79=\"],\n  [\"940\u00f76=\", \"239\u00f77=\"],\n  [\"202\u00f73=\", \"761\u00f72=\"],\n  [\"834\u00f72=\", \"988\u00f79=\"],\n  [\"507\u00f79=\", \"589\u00f73=\"],\n  [\"379\u00f78=\", \"972\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-divided-by-one-digit expression in the\n# document with its updated value, per the commit's regenerated\n# worksheet numbers. Every \"find\" string is a unique, literal run of\n# text (e.g. \"717\u00f79=\"), so Find/Replace on the whole-document range\n# swaps just the text while leaving the surrounding run formatting\n# (font, size, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"717\u00f79=\", \"147\u00f76=\"),\n    @(\"632\u00f77=\", \"914\u00f75=\"),\n    @(\"169\u00f74=\", \"873\u00f77=\"),\n    @(\"852\u00f79=\", \"153\u00f72=\"),\n    @(\"473\u00f78=\", \"562\u00f72=\"),\n    @(\"531\u00f74=\", \"235\u00f74=\"),\n    @(\"284\u00f72=\", \"795\u00f73=\"),\n    @(\"431\u00f79=\", \"503\u00f72=\"),\n    @(\"946\u00f78=\", \"891\u00f77=\"),\n    @(\"804\u00f79=\", \"267\u00f76=\"),\n    @(\"976\u00f76=\", \"568\u00f77=\"),\n    @(\"951\u00f78=\", \"676\u00f73=\"),\n    @(\"169\u00f72=\", \"405\u00f78=\"),\n    @(\"904\u00f78=\", \"583\u00f75=\"),\n    @(\"694\u00f72=\", \"869\u00f73=\"),\n    @(\"863\u00f75=\", \"934\u00f72=\"),\n    @(\"724\u00f72=\", \"542\u00f75=\"),\n    @(\"116\u00f78=\", \"814\u00f74=\"),\n    @(\"861\u00f73=\", \"763\u00f77=\"),\n    @(\"217\u00f73=\", \"840\u00f79=\"),\n    @(\"940\u00f76=\", \"239\u00f77=\"),\n    @(\"202\u00f73=\", \"761\u00f72=\"),\n    @(\"834\u00f72=\", \"988\u00f79=\"),\n    @(\"507\u00f79=\", \"589\u00f73=\"),\n    @(\"379\u00f78=\", \"972\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
